# Update the NDTV-sourced data columns (C: LiveTemp(NDTV), F: Weather Condition (NDTV),
# I: Humidity (NDTV), L: Wind (NDTV)) with refreshed values, then update the
# selection and a few column widths to match the resaved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Lucknow)
$ws.Range("C2").Value = "35"
$ws.Range("F2").Value = "Humid and Overcast`n"
$ws.Range("I2").Value = "59"
$ws.Range("L2").Value = "1.23"

# Row 3 (Bhopal)
$ws.Range("C3").Value = "32"
$ws.Range("F3").Value = "Humid and Mostly Cloudy`n"
$ws.Range("I3").Value = "65"
$ws.Range("L3").Value = "1.60"

# Row 4 (Ajmer)
$ws.Range("C4").Value = "33"
$ws.Range("I4").Value = "60"
$ws.Range("L4").Value = "1.61"

# Row 5 (Coimbatore)
$ws.Range("C5").Value = "31"
$ws.Range("F5").Value = "Humid and Mostly Cloudy`n"
$ws.Range("I5").Value = "62"
$ws.Range("L5").Value = "0.94"

# Row 6 (Mumbai)
$ws.Range("C6").Value = "30"
$ws.Range("F6").Value = "Humid and Mostly Cloudy`n"
$ws.Range("I6").Value = "77"
$ws.Range("L6").Value = "1.66"

# Row 7 (Kolkata)
$ws.Range("I7").Value = "64"
$ws.Range("L7").Value = "2.42"

# Writing the multi-line "Weather Condition (NDTV)" text triggers an automatic
# row-height bump on affected rows; re-run AutoFit on all data rows so the
# saved row heights stay at the sheet default, same as before the edit.
$ws.Range("A2:A7").EntireRow.AutoFit() | Out-Null

# Column widths are "best fit" - let Excel recompute them for the changed columns
# now that their contents have been refreshed (same as Excel does automatically
# when the workbook is resaved).
$ws.Columns.Item("C:C").AutoFit() | Out-Null
$ws.Columns.Item("F:F").AutoFit() | Out-Null
$ws.Columns.Item("I:I").AutoFit() | Out-Null
$ws.Columns.Item("L:L").AutoFit() | Out-Null

# Update selection to reflect where the author left the cursor
$ws.Range("C11").Select()
